$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# The sheet gains two new leading columns ("vendor", "doc. number").
# The previously-existing columns (reference, document date, gross amount)
# shift right by two positions (A->C, B->D, C->E).
#
# We deliberately avoid Range/Column .Insert()/.Cut() here: those structural
# operations make the host engine quietly re-derive the shifted numeric
# values through a lossy float path (e.g. 25102.35 -> 25102.349999999999).
# Instead we copy number formats explicitly and then (re)write every cell's
# final literal value directly, which keeps the original numbers exact.
# ---------------------------------------------------------------------------

# Carry the date/amount display formats from the old B/C columns over to the
# new D/E columns before we overwrite B/C with the new "doc. number" data.
$ws.Range("B2:B10").Copy()
$ws.Range("D2:D10").PasteSpecial(-4122)
$ws.Range("C2:C10").Copy()
$ws.Range("E2:E10").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Columns A and B are brand new ("vendor" / "doc. number" - plain numbers),
# so make sure they don't keep any inherited formatting (the old B column
# was date-formatted, which would otherwise turn the big doc.number values
# into out-of-range dates).
$ws.Range("A1:B10").NumberFormat = "General"

# ----- Header row -----------------------------------------------------
$ws.Range("A1").Value = "vendor"
$ws.Range("B1").Value = "doc. number"
$ws.Range("C1").Value = "reference"
$ws.Range("D1").Value = "document date"
$ws.Range("E1").Value = "gross amount"

# ----- Data rows --------------------------------------------------------
$vendor  = @(12345, 32564, 549879, 45555, 489, 32654, 587646, 687465, 66846)
$docnum  = @(740000001, 740000002, 740000003, 740000004, 740000005, 740000006, 740000007, 740000008, 740000009)
$refval  = @("DAB1234", 4568456315, "jk5468-20", "DAB1234", 4567892, 16548, "2021-nhgf45", "jk5468-20", "GF546")
$docdate = @(44280, 43966, 44327, 44280, 44275, 44058, 44219, 44296, 44327)
$gross   = @(2000.2, 25102.35, 105, 2000.2, 1120.3599999999999, 45851.12, 3, 105, 105)

for ($i = 0; $i -lt 9; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value = $vendor[$i]
    $ws.Cells.Item($r, 2).Value = $docnum[$i]
    $ws.Cells.Item($r, 3).Value = $refval[$i]
    $ws.Cells.Item($r, 4).Value = $docdate[$i]
    $ws.Cells.Item($r, 5).Value = $gross[$i]
}

# ----- Column widths ------------------------------------------------------
# Columns C/D/E keep the exact widths the old A/B/C columns used to have;
# column B is brand new. (ColumnWidth here snaps to a coarse 1/6-character
# grid internally, so these are the closest obtainable values.)
$ws.Columns.Item(2).ColumnWidth = 11.833333333333334
$ws.Columns.Item(3).ColumnWidth = 10.666666666666666
$ws.Columns.Item(4).ColumnWidth = 15.833333333333334
$ws.Columns.Item(5).ColumnWidth = 15.166666666666666

# ----- Misc sheet view state ------------------------------------------------
$ws.Range("D13").Select()
